# ------------------------------------------------------------------------
# Adds a new "Player Info" sheet (as the first sheet of the workbook) that
# holds basic player metadata, and updates the existing "ODI Batting" /
# "ODI Bowling" sheets so that the old MATCH_CARD_LINK (a full scorecard
# URL) column becomes a much shorter MATCH_CODE column (just the numeric
# match code extracted from that URL).
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE ---------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "'4675"

# --- 2. Update "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE ---------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "'4675"

# --- 3. Create the new "Player Info" sheet as the first sheet -------------
# Cloning an existing sheet (instead of Worksheets.Add()) lets the new
# sheet inherit the same header styling (bold / bordered / centered) that
# is already used on the other sheets, without introducing new styles.
$batting.Copy($batting)
$playerInfo = $wb.Worksheets.Item("ODI Batting (2)")
$playerInfo.Name = "Player Info"

# Drop the extra columns (E:J) that came along with the clone - "Player
# Info" only needs 4 columns.
$playerInfo.Range("E1:J2").Clear()

# Headers
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row
$playerInfo.Range("A2").Value = "'6536"
$playerInfo.Range("B2").Value = "Noor Ahmad Lakanwal"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Left Arm Wrist Spin (Chinaman)"
